$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update District column (G) values from "Mysore" (and variants like
# "Pereyapatna" / trailing-space "Mysore ") to the official name
# "Mysuru (Mysore)" per the website's corrected naming.
for ($r = 4; $r -le 53; $r++) {
    $ws.Cells.Item($r, 7).Value = "Mysuru (Mysore)"
}
